# Fruta / hortaliza, semanal
# Insert a new weekly observation row at row 20 (pushing existing rows 20-37
# down to 21-38) and populate it with the new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 20; Excel shifts rows 20:37 to
# rows 21:38 and picks up formatting (e.g. the date style on column D)
# from the row above, matching the existing data rows.
$ws.Rows.Item(20).Insert()

$ws.Cells.Item(20, 1).Value2 = 6
$ws.Cells.Item(20, 2).Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(20, 3).Value2 = "Metropolitana"
$ws.Cells.Item(20, 4).Value2 = 44790
$ws.Cells.Item(20, 5).Value2 = 13
$ws.Cells.Item(20, 6).Value2 = 100112035
$ws.Cells.Item(20, 7).Value2 = "Bruselas (repollito)"
$ws.Cells.Item(20, 8).Value2 = "Sin especificar"
$ws.Cells.Item(20, 9).Value2 = "Primera"
$ws.Cells.Item(20, 10).Value2 = 500
$ws.Cells.Item(20, 11).Value2 = 15000
$ws.Cells.Item(20, 12).Value2 = 16000
$ws.Cells.Item(20, 13).Value2 = 15540
$ws.Cells.Item(20, 14).Value2 = "`$/malla 15 kilos"
$ws.Cells.Item(20, 15).Value2 = "Provincia de Quillota"
$ws.Cells.Item(20, 16).Value2 = 1036
$ws.Cells.Item(20, 17).Value2 = 15
$ws.Cells.Item(20, 18).Value2 = "Hortaliza"
